$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (for the two newest fiscal quarters)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formatting (date / number styles) from old column D (now column F)
# across into the two newly inserted columns D and E
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the values for the two new columns, plus a handful of corrected
# historical figures that landed in column H after the insert

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 728200
$ws.Range("E8").Value = 309700
$ws.Range("H8").Value = 704100
$ws.Range("D9").Value = 585000
$ws.Range("E9").Value = 249100
$ws.Range("H9").Value = 553000
$ws.Range("D10").Value = 143200
$ws.Range("E10").Value = 60600
$ws.Range("H10").Value = 151100
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 745600
$ws.Range("E17").Value = 298600
$ws.Range("H17").Value = 716800
$ws.Range("D18").Value = -17400
$ws.Range("E18").Value = 11100
$ws.Range("H18").Value = -12700
$ws.Range("D20").Value = -2900
$ws.Range("E20").Value = -900
$ws.Range("H20").Value = -3700
$ws.Range("D21").Value = 34800
$ws.Range("E21").Value = 24500
$ws.Range("H21").Value = 62300
$ws.Range("D22").Value = 4800
$ws.Range("E22").Value = 900
$ws.Range("H22").Value = 2700
$ws.Range("D23").Value = -25100
$ws.Range("E23").Value = 9300
$ws.Range("H23").Value = -19100
$ws.Range("D24").Value = -6100
$ws.Range("E24").Value = 1400
$ws.Range("H24").Value = -2800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -19000
$ws.Range("E26").Value = 7900
$ws.Range("H26").Value = -16300
$ws.Range("D27").Value = -19000
$ws.Range("E27").Value = 7900
$ws.Range("H27").Value = -16300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -144000
$ws.Range("E29").Value = "NA"
$ws.Range("H29").Value = -25500
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 2900
$ws.Range("E32").Value = 900
$ws.Range("H32").Value = 3700
$ws.Range("D33").Value = -163000
$ws.Range("E33").Value = 7900
$ws.Range("H33").Value = -41800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -163000
$ws.Range("E35").Value = 7900
$ws.Range("H35").Value = -41800
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 173800
$ws.Range("E41").Value = 39400
$ws.Range("H41").Value = 40100
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 39600
$ws.Range("E43").Value = 116800
$ws.Range("H43").Value = 56200
$ws.Range("D44").Value = 139500
$ws.Range("E44").Value = 266600
$ws.Range("H44").Value = 117000
$ws.Range("D45").Value = 18000
$ws.Range("E45").Value = 22500
$ws.Range("H45").Value = 208800
$ws.Range("D46").Value = 370900
$ws.Range("E46").Value = 445300
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 296100
$ws.Range("E48").Value = 398900
$ws.Range("H48").Value = 310600
$ws.Range("D49").Value = 32000
$ws.Range("E49").Value = 55400
$ws.Range("H49").Value = 99000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 9200
$ws.Range("E52").Value = 11100
$ws.Range("H52").Value = 150600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 708200
$ws.Range("E54").Value = 910700
$ws.Range("H54").Value = 890100
$ws.Range("D57").Value = 99500
$ws.Range("E57").Value = 127800
$ws.Range("H57").Value = 93600
$ws.Range("D58").Value = 28700
$ws.Range("E58").Value = 200
$ws.Range("D59").Value = 42400
$ws.Range("E59").Value = 43000
$ws.Range("H59").Value = 75300
$ws.Range("D60").Value = 170600
$ws.Range("E60").Value = 171000
$ws.Range("D61").Value = 70600
$ws.Range("E61").Value = 94600
$ws.Range("H61").Value = 85000
$ws.Range("D62").Value = 76000
$ws.Range("E62").Value = 88300
$ws.Range("H62").Value = 104100
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 317200
$ws.Range("E66").Value = 353900
$ws.Range("H66").Value = 340100
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -186300
$ws.Range("E72").Value = -15300
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 391000
$ws.Range("E76").Value = 556800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -163000
$ws.Range("E81").Value = 7900
$ws.Range("H81").Value = -41800
$ws.Range("D83").Value = 55100
$ws.Range("E83").Value = 14300
$ws.Range("H83").Value = 78700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 62500
$ws.Range("E89").Value = -6800
$ws.Range("H89").Value = 62900
$ws.Range("D91").Value = -35300
$ws.Range("E91").Value = -6900
$ws.Range("H91").Value = -44800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 60600
$ws.Range("E94").Value = -6900
$ws.Range("H94").Value = -80500
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 13300
$ws.Range("E100").Value = 25300
$ws.Range("H100").Value = 24400
$ws.Range("D101").Value = -1600
$ws.Range("E101").Value = -700
$ws.Range("H101").Value = 1600
$ws.Range("D102").Value = 134800
$ws.Range("E102").Value = 10900
$ws.Range("H102").Value = 8400
